# Add season-record columns (Wins / Losses / Ties) to the player table.
# New columns AD:AF are appended after the existing AC column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the existing header formatting (bold, border,
#     centered) from A1 onto the new header cells, then set their text. ---
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-50): every player gets their team's season record. ---
$lastRow = 50
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 88   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
